$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit adds a small "OrigenProducto" mini-table (M:O columns) next to
# each of the four "Producto" tables on the sheet (rows 12-13, 33-34, 53-54,
# 73-74). The header cell (row 12/33/53/73, column M) reuses the same visual
# style as the existing table headers (e.g. E12 "Producto"), and the three
# detail cells (row 13/34/54/74, columns M,N,O) reuse the same "key column"
# style as the existing blue/bold-italic header cells (e.g. D4, K13).
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122
$xlNone = -4142
$xlContinuous = 1
$xlBlack = 1

# Build a reusable "key" style (bold italic white on blue, thin left+right
# border only) once, in a scratch cell far away from the used range, then
# fan it out to every target cell. Doing it this way keeps every target
# cell pointing at the exact same cell style (same <xf>) instead of each
# accumulating its own near-duplicate style.
$tmpl = $ws.Range("Z100")
$ws.Range("D4").Copy()
$tmpl.PasteSpecial($xlPasteFormats)
$tmpl.Borders.LineStyle = $xlNone
$tmpl.Borders.Item(7).ColorIndex = $xlBlack
$tmpl.Borders.Item(7).LineStyle = $xlContinuous
$tmpl.Borders.Item(10).ColorIndex = $xlBlack
$tmpl.Borders.Item(10).LineStyle = $xlContinuous

$tmpl.Copy()
foreach ($row in @(13, 34, 54, 74)) {
    $ws.Range("M$row").PasteSpecial($xlPasteFormats)
    $ws.Range("N$row").PasteSpecial($xlPasteFormats)
    $ws.Range("O$row").PasteSpecial($xlPasteFormats)
}

# Header cell (row above), reusing the existing "Producto" header style.
foreach ($row in @(12, 33, 53, 73)) {
    $ws.Range("E$row").Copy()
    $ws.Range("M$row").PasteSpecial($xlPasteFormats)
    $ws.Range("M$row").Value = "OrigenProducto"
}

# Fill in the detail row values.
foreach ($row in @(13, 34, 54, 74)) {
    $ws.Range("M$row").Value = "idProducto"
    $ws.Range("N$row").Value = "idProveedor"
    $ws.Range("O$row").Value = "RIF"
}

$tmpl.Clear()

$ws.Range("M12:O13").Select()
